# Apply LOQ4087.xlsx content realignment (rows 13-26) + new text blocks
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear cells that must become empty after the realignment ---
$ws.Range("A13").ClearContents()
$ws.Range("B18").ClearContents()
$ws.Range("C18").ClearContents()
$ws.Range("B23").ClearContents()
$ws.Range("C23").ClearContents()

# --- Set cell values (A/B/C) for every populated row, 1 through 26 ---
$ws.Range("B1").Value = "Ementa atual:"
$ws.Range("C1").Value = "Ementa modificada (dados modificados em vermelho):"

$ws.Range("B2").Value = "LOQ4087"
$ws.Range("C2").Value = "LOQ4087"

$ws.Range("A3").Value = "Nome:"
$ws.Range("B3").Value = " Termodinâmica Química Aplicada I"
$ws.Range("C3").Value = " Termodinâmica Química Aplicada I"

$ws.Range("A4").Value = "Name:"
$ws.Range("B4").Value = "Applied Chemical Thermodynamics I"
$ws.Range("C4").Value = "Applied Chemical Thermodynamics I"

$ws.Range("A5").Value = "Créditos-aula:"
$ws.Range("B5").Value = "4"
$ws.Range("C5").Value = "4"

$ws.Range("A6").Value = "Créditos-trabalho"
$ws.Range("B6").Value = "0"
$ws.Range("C6").Value = "0"

$ws.Range("A7").Value = "Carga horária:"
$ws.Range("B7").Value = "60 h"
$ws.Range("C7").Value = "60 h"

$ws.Range("A8").Value = "Ativação:"
$ws.Range("B8").Value = "01/01/2020"
$ws.Range("C8").Value = "01/01/2020"

$ws.Range("A9").Value = "Semestre ideal:"
$ws.Range("B9").Value = "EB-4,EQD-4,EQN-5"
$ws.Range("C9").Value = "EB-4,EQD-4,EQN-5"

$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "Gerais Levar os estudantes a compreenderem que a Termodinâmica, juntamente com fenomenos de transporte constituem-se na base fundamental para o entendimento e solução de grande parte dos problemas de Eng. Química. 2) Específicos Ao final do curso os estudantes deverão: - Compreender os aspectos mássicos, energéticos e entrópicos, envolvendo sistemas termodinâmicos abertos e fechados; - Dominar e ser capaz de fazer predições de propriedades termodinâmicas, usando equações cúbicas de estado e relações termodinâmicas; - Desenvolver uma metodologia para poder solucionar os problemas de engenharia, nos aspectos termodinâmicos; Dominar o uso de tabelas de propriedades termodinâmicas;"
$ws.Range("C10").Value = "Gerais Levar os estudantes a compreenderem que a Termodinâmica, juntamente com fenomenos de transporte constituem-se na base fundamental para o entendimento e solução de grande parte dos problemas de Eng. Química. 2) Específicos Ao final do curso os estudantes deverão: - Compreender os aspectos mássicos, energéticos e entrópicos, envolvendo sistemas termodinâmicos abertos e fechados; - Dominar e ser capaz de fazer predições de propriedades termodinâmicas, usando equações cúbicas de estado e relações termodinâmicas; - Desenvolver uma metodologia para poder solucionar os problemas de engenharia, nos aspectos termodinâmicos; Dominar o uso de tabelas de propriedades termodinâmicas;"

$ws.Range("A11").Value = "Objectives:"
$ws.Range("B11").Value = "General objectiveStudents will be led to understand that thermodynamics, together with transport phenomena, are the fundamental basis for understanding and solving most of the problems of Chemical Engineering.Specific objectivesAt the end of the course, students will understand the mass, energetic and entropic aspects, involving open and closed thermodynamic systems, to master and being able to make predictions of properties using cubic equations of state and thermodynamic relations, to develop a methodology for solving engineering problems in thermodynamic aspects, to master the use of tables of thermodynamic properties."
$ws.Range("C11").Value = "General objectiveStudents will be led to understand that thermodynamics, together with transport phenomena, are the fundamental basis for understanding and solving most of the problems of Chemical Engineering.Specific objectivesAt the end of the course, students will understand the mass, energetic and entropic aspects, involving open and closed thermodynamic systems, to master and being able to make predictions of properties using cubic equations of state and thermodynamic relations, to develop a methodology for solving engineering problems in thermodynamic aspects, to master the use of tables of thermodynamic properties."

$ws.Range("A12").Value = "Docentes responsáveis:"

$ws.Range("B13").Value = "8554681 - Pedro Felipe Arce Castillo"
$ws.Range("C13").Value = "8554681 - Pedro Felipe Arce Castillo"

$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = "Conceitos e definições básicas. Propriedades de substâncias puras. Primeira lei da termodinâmica. Entropia e segunda lei da termodinâmica. Termodinâmica de processos de escoamento. Produção de potência a partir de calor. Refrigeração e liquefação."
$ws.Range("C14").Value = "Conceitos e definições básicas. Propriedades de substâncias puras. Primeira lei da termodinâmica. Entropia e segunda lei da termodinâmica. Termodinâmica de processos de escoamento. Produção de potência a partir de calor. Refrigeração e liquefação."

$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("B15").Value = "Concepts and basic definitions. Properties of pure substances. First law of thermodynamics. Entropy and second law of thermodynamics. Thermodynamics of flow processes. Production of power from heat. Refrigeration and liquefaction."
$ws.Range("C15").Value = "Concepts and basic definitions. Properties of pure substances. First law of thermodynamics. Entropy and second law of thermodynamics. Thermodynamics of flow processes. Production of power from heat. Refrigeration and liquefaction."

$ws.Range("A16").Value = "Programa:"
$ws.Range("B16").Value = "1 - Conceitos e definições básicas.1.1-Termodinâmica.1.2-Sistema termodinâmico.1.3-Substância pura, Estado, Fase.1.4-Função de estado, função caminho.1.5-Propriedades termodinâmicas: Pressão, Temperatura, Volume.1.6-Equilíbrio termodinâmico.1.7-Processos termodinâmicos.1.8-Processos reversíveis e irreversíveis.1.9-Ciclos termodinâmicos: Ciclo de potência a vapor e Ciclo de refrigeração.1.10-Calor e Trabalho.1.11-Balanço de uma propriedade termodinâmica.2 - Propriedades de substâncias puras.2.1- Substância pura. Diagrama PTv2.2-Pressão de vapor/Pressão de saturação2.3-Ponto crítico e Ponto triplo2.4-Critérios para determinar os estados termodinâmicos nas regiões de saturação, líquido comprimido e vapor superaquecido utilizando tabelas de propriedades termodinâmicas.2.5-Diagramas termodinâmicos vs Tabelas termodinâmicas2.6-Equação de estado para fluido puro2.7-Exercícios.3 - Primeira lei da termodinâmica.3.1-Definição. Conservação da energia3.2-Balanço de massa e a equação da continuidade3.3-Energia: energia interna, potencial e cinética3.4-Equação geral da primeira Lei da termodinâmica: Balanço de energia3.5-Propriedades termodinâmicas: energia interna e entalpia nas regiões de saturação, líquido comprimido e vapor superaquecido.3.6-Calor específico a volume e pressão constante3.7-Grandezas: Calor e Trabalho. Convenção de sinais.3.8-Exercícios.4 - Entropia e segunda lei da termodinâmica.4.1-Definição. Geração de entropia4.2-Primeira lei versus Segunda lei da termodinâmica4.3-Análise da segunda lei da termodinâmica4.4-Máquinas térmicas de potência. Eficiência4.5-Enunciados da segunda lei da termodinâmica: Kelvin-Planck e Clausius.4.6-Processos reversíveis e irreversíveis4.7-Irreversibilidades.4.8-Quantificação da entropia: Desigualdade de Clausius4.9-Propriedade termodinâmica: entropia nas regiões de saturação, líquido comprimido e vapor superaquecido.4.9-Entropia gerada4.10-Equação geral da segunda Lei da termodinâmica: Balanço de entropia.4.11-Eficiência de uma máquina térmica.4.12-Processos e eficiência isentrópica: Turbina, compressor, bomba4.13-Exercícios.5 - Termodinâmica de processos de escoamento.5.1-Escoamento de fluidos. Balanço de energia mecânica. Equação de Bernoulli5.2-Expansão de fluidos: turbinas e válvulas de expansão.5.3-Compressão de fluidos: compressores e bombas.5.4-Troca de calor: Trocadores de calor.5.5-Aumento da energia cinética: bocal5.6-Exercícios.6 - Produção de potência a partir de calor.6.1-Ciclos de potência a vapor.6.2-Ciclo Carnot.6.3-Ciclo Rankine: ideal e real. Eficiência6.4-Tipos de Ciclos Rankine: com Reaquecimento e Regenerativo.6.5-Exercícios.6.6-Ciclo Brayton: ideal e real. Eficiência6.7-Ciclo combinado: Brayton-Rankine.6.8-Exercícios.7 - Refrigeração e liquefação.7.1-Ciclo de refrigeração por compressão de vapor: ideal e real.7.2-Configurações dos ciclos de refrigeração: em dois estágios e cascata.7.3-Fluidos de trabalho utilizados.7.4-Aplicações domésticas e na engenharia.7.5-Liquefação de gases.7.6-Exercícios."
$ws.Range("C16").Value = "1 - Conceitos e definições básicas.1.1-Termodinâmica.1.2-Sistema termodinâmico.1.3-Substância pura, Estado, Fase.1.4-Função de estado, função caminho.1.5-Propriedades termodinâmicas: Pressão, Temperatura, Volume.1.6-Equilíbrio termodinâmico.1.7-Processos termodinâmicos.1.8-Processos reversíveis e irreversíveis.1.9-Ciclos termodinâmicos: Ciclo de potência a vapor e Ciclo de refrigeração.1.10-Calor e Trabalho.1.11-Balanço de uma propriedade termodinâmica.2 - Propriedades de substâncias puras.2.1- Substância pura. Diagrama PTv2.2-Pressão de vapor/Pressão de saturação2.3-Ponto crítico e Ponto triplo2.4-Critérios para determinar os estados termodinâmicos nas regiões de saturação, líquido comprimido e vapor superaquecido utilizando tabelas de propriedades termodinâmicas.2.5-Diagramas termodinâmicos vs Tabelas termodinâmicas2.6-Equação de estado para fluido puro2.7-Exercícios.3 - Primeira lei da termodinâmica.3.1-Definição. Conservação da energia3.2-Balanço de massa e a equação da continuidade3.3-Energia: energia interna, potencial e cinética3.4-Equação geral da primeira Lei da termodinâmica: Balanço de energia3.5-Propriedades termodinâmicas: energia interna e entalpia nas regiões de saturação, líquido comprimido e vapor superaquecido.3.6-Calor específico a volume e pressão constante3.7-Grandezas: Calor e Trabalho. Convenção de sinais.3.8-Exercícios.4 - Entropia e segunda lei da termodinâmica.4.1-Definição. Geração de entropia4.2-Primeira lei versus Segunda lei da termodinâmica4.3-Análise da segunda lei da termodinâmica4.4-Máquinas térmicas de potência. Eficiência4.5-Enunciados da segunda lei da termodinâmica: Kelvin-Planck e Clausius.4.6-Processos reversíveis e irreversíveis4.7-Irreversibilidades.4.8-Quantificação da entropia: Desigualdade de Clausius4.9-Propriedade termodinâmica: entropia nas regiões de saturação, líquido comprimido e vapor superaquecido.4.9-Entropia gerada4.10-Equação geral da segunda Lei da termodinâmica: Balanço de entropia.4.11-Eficiência de uma máquina térmica.4.12-Processos e eficiência isentrópica: Turbina, compressor, bomba4.13-Exercícios.5 - Termodinâmica de processos de escoamento.5.1-Escoamento de fluidos. Balanço de energia mecânica. Equação de Bernoulli5.2-Expansão de fluidos: turbinas e válvulas de expansão.5.3-Compressão de fluidos: compressores e bombas.5.4-Troca de calor: Trocadores de calor.5.5-Aumento da energia cinética: bocal5.6-Exercícios.6 - Produção de potência a partir de calor.6.1-Ciclos de potência a vapor.6.2-Ciclo Carnot.6.3-Ciclo Rankine: ideal e real. Eficiência6.4-Tipos de Ciclos Rankine: com Reaquecimento e Regenerativo.6.5-Exercícios.6.6-Ciclo Brayton: ideal e real. Eficiência6.7-Ciclo combinado: Brayton-Rankine.6.8-Exercícios.7 - Refrigeração e liquefação.7.1-Ciclo de refrigeração por compressão de vapor: ideal e real.7.2-Configurações dos ciclos de refrigeração: em dois estágios e cascata.7.3-Fluidos de trabalho utilizados.7.4-Aplicações domésticas e na engenharia.7.5-Liquefação de gases.7.6-Exercícios."

$ws.Range("A17").Value = "Syllabus:"
$ws.Range("B17").Value = "1 - Basic concepts and definitions.1.1-Thermodynamics.1.2-Thermodynamic system.1.3-Pure substance, State, Phase.1.4-State function, path function.1.5-Thermodynamic Properties: Pressure, Temperature, Volume.1.6-Thermodynamic equilibrium.1.7-Thermodynamic processes.1.8-Reversible and irreversible processes.1.9-Thermodynamic cycles: Steam power cycle and Refrigeration cycle.1.10-Heat and Work.1.11-Balance of a thermodynamic property.2 - Properties of pure substances.2.1- Pure substance. PTv diagram2.2-Vapor Pressure / Saturation Pressure2.3-Critical point and Triple point2.4-Criteria for determining the thermodynamic states in the regions of saturation, compressed liquid and superheated steam by using tables of thermodynamic properties.2.5-Thermodynamic diagrams vs Thermodynamic tables2.6-Equation of State for pure fluid2.7-Exercises.3 - First law of thermodynamics. 3.1-Definition. Energy Conservation 3.2-Mass balance and the continuity equation 3.3-Energy: internal, potential and kinetic energy 3.4-General equation of the first law of thermodynamics: Energy balance 3.5-Thermodynamic properties: internal energy and enthalpy in the regions of saturation, compressed liquid and superheated steam. 3.6-Specific heat at constant volume and pressure 3.7-Greatness: Heat and Work. Sign Convention. 3.8-Exercises.4 - Entropy and second law of thermodynamics.4.1-Definition. Generation of entropy4.2-First law versus Second law of thermodynamics4.3-Analysis of the second law of thermodynamics4.4-Thermal power machines. Efficiency4.5-Statements of the second law of thermodynamics: Kelvin-Planck and Clausius.4.6-Reversible and irreversible processes4.7-Irreversibilities.4.8-Entropy Quantification: Clausius Inequality4.9-Thermodynamic properties: entropy in the regions of saturation, compressed liquid and superheated steam.4.9-Entropy generated4.10-General equation of the second law of thermodynamics: Balance of entropy.4.11-Efficiency of a thermal machine.4.12-Isentropic processes and efficiency: Turbine, compressor, pump4.13-Exercises.5 - Thermodynamics of flow processes.5.1-Flow of fluids. Balance of mechanical energy. Bernoulli equation5.2-Expansion of fluids: turbines and expansion valves.5.3-Compression of fluids: Compressors and pumps.5.4-Heat exchange: heat exchangers.5.5-Increase of kinetic energy: nozzle5.6-Exercises.6 - Production of power from heat.6.1-Steam power cycles.6.2-Carnot Cycle.6.3-Rankine Cycle: ideal and real. Efficiency6.4-Types of Rankine Cycles: with Reheating and Regenerative.6.5-Exercises.6.6-Brayton cycle: ideal and real. Efficiency6.7-Combined cycle: Brayton-Rankine.6.8-Exercises.7 - Refrigeration and liquefaction.7.1-Cycle of vapor-compression refrigeration: ideal and real.7.2-Configurations of the refrigeration cycles: two stages and cascade.7.3-Working fluids used.7.4-Domestic and engineering applications.7.5-Liquefaction of gases.7.6-Exercises."
$ws.Range("C17").Value = "1 - Basic concepts and definitions.1.1-Thermodynamics.1.2-Thermodynamic system.1.3-Pure substance, State, Phase.1.4-State function, path function.1.5-Thermodynamic Properties: Pressure, Temperature, Volume.1.6-Thermodynamic equilibrium.1.7-Thermodynamic processes.1.8-Reversible and irreversible processes.1.9-Thermodynamic cycles: Steam power cycle and Refrigeration cycle.1.10-Heat and Work.1.11-Balance of a thermodynamic property.2 - Properties of pure substances.2.1- Pure substance. PTv diagram2.2-Vapor Pressure / Saturation Pressure2.3-Critical point and Triple point2.4-Criteria for determining the thermodynamic states in the regions of saturation, compressed liquid and superheated steam by using tables of thermodynamic properties.2.5-Thermodynamic diagrams vs Thermodynamic tables2.6-Equation of State for pure fluid2.7-Exercises.3 - First law of thermodynamics. 3.1-Definition. Energy Conservation 3.2-Mass balance and the continuity equation 3.3-Energy: internal, potential and kinetic energy 3.4-General equation of the first law of thermodynamics: Energy balance 3.5-Thermodynamic properties: internal energy and enthalpy in the regions of saturation, compressed liquid and superheated steam. 3.6-Specific heat at constant volume and pressure 3.7-Greatness: Heat and Work. Sign Convention. 3.8-Exercises.4 - Entropy and second law of thermodynamics.4.1-Definition. Generation of entropy4.2-First law versus Second law of thermodynamics4.3-Analysis of the second law of thermodynamics4.4-Thermal power machines. Efficiency4.5-Statements of the second law of thermodynamics: Kelvin-Planck and Clausius.4.6-Reversible and irreversible processes4.7-Irreversibilities.4.8-Entropy Quantification: Clausius Inequality4.9-Thermodynamic properties: entropy in the regions of saturation, compressed liquid and superheated steam.4.9-Entropy generated4.10-General equation of the second law of thermodynamics: Balance of entropy.4.11-Efficiency of a thermal machine.4.12-Isentropic processes and efficiency: Turbine, compressor, pump4.13-Exercises.5 - Thermodynamics of flow processes.5.1-Flow of fluids. Balance of mechanical energy. Bernoulli equation5.2-Expansion of fluids: turbines and expansion valves.5.3-Compression of fluids: Compressors and pumps.5.4-Heat exchange: heat exchangers.5.5-Increase of kinetic energy: nozzle5.6-Exercises.6 - Production of power from heat.6.1-Steam power cycles.6.2-Carnot Cycle.6.3-Rankine Cycle: ideal and real. Efficiency6.4-Types of Rankine Cycles: with Reheating and Regenerative.6.5-Exercises.6.6-Brayton cycle: ideal and real. Efficiency6.7-Combined cycle: Brayton-Rankine.6.8-Exercises.7 - Refrigeration and liquefaction.7.1-Cycle of vapor-compression refrigeration: ideal and real.7.2-Configurations of the refrigeration cycles: two stages and cascade.7.3-Working fluids used.7.4-Domestic and engineering applications.7.5-Liquefaction of gases.7.6-Exercises."

$ws.Range("A18").Value = "Avaliação:"

$ws.Range("A19").Value = "Método:"
$ws.Range("B19").Value = "2 provas escritas"
$ws.Range("C19").Value = "2 provas escritas"

$ws.Range("A20").Value = "Critério:"
$ws.Range("B20").Value = "Serão avaliados os conteúdos discutidos em sala e constantes da ementa do curso. A média da disciplina será a média aritmética das duas provas."
$ws.Range("C20").Value = "Serão avaliados os conteúdos discutidos em sala e constantes da ementa do curso. A média da disciplina será a média aritmética das duas provas."

$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Range("B21").Value = "prova escrita com conteúdo de todo o semestre"
$ws.Range("C21").Value = "prova escrita com conteúdo de todo o semestre"

$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = "KORETSKY, M. D. Termodinâmica para Engenharia Química. 1 ed. LTC Editora, 2007. `nMORAN, M. J.; SHAPIRO, H. N. Princípios de Termodinâmica para Engenharia. 1 ed. LTC Editora, 2009. `nSANDLER, S. I., Chemical and Engineering Thermodynamics, 3rd ed., John Wiley & Sons, 1999 `nSMITH, J.M.; VAN NESS, H.C.; Abott, M. M. Introdução à Termodinâmica da Engenharia Química. 7ª ed. LTC editora, 2007. `nTERRON, L. R. Termodinâmica Química Aplicada. 1 ed. Editora Manole Ltda, 2009. `nVAN WILEN, J. Sonntag, Richard. E. Fundamentos da Termodinâmica Clássica. 6 ed. 2004"
$ws.Range("C22").Value = "KORETSKY, M. D. Termodinâmica para Engenharia Química. 1 ed. LTC Editora, 2007. `nMORAN, M. J.; SHAPIRO, H. N. Princípios de Termodinâmica para Engenharia. 1 ed. LTC Editora, 2009. `nSANDLER, S. I., Chemical and Engineering Thermodynamics, 3rd ed., John Wiley & Sons, 1999 `nSMITH, J.M.; VAN NESS, H.C.; Abott, M. M. Introdução à Termodinâmica da Engenharia Química. 7ª ed. LTC editora, 2007. `nTERRON, L. R. Termodinâmica Química Aplicada. 1 ed. Editora Manole Ltda, 2009. `nVAN WILEN, J. Sonntag, Richard. E. Fundamentos da Termodinâmica Clássica. 6 ed. 2004"

$ws.Range("A23").Value = "Requisitos:"

$ws.Range("B24").Value = "LOB1004 -  Cálculo II  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOB1004 -  Cálculo II  (Requisito fraco)`n"

$ws.Range("B25").Value = "LOB1019 -  Física II  (Requisito fraco)`n"
$ws.Range("C25").Value = "LOB1019 -  Física II  (Requisito fraco)`n"

$ws.Range("B26").Value = "LOQ4053 -  Balanços de Massa e Energia  (Requisito fraco)`n"
$ws.Range("C26").Value = "LOQ4053 -  Balanços de Massa e Energia  (Requisito fraco)`n"

# --- Row heights: match the target layout exactly ---
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 120
$ws.Rows.Item(23).AutoFit()
$ws.Rows.Item(24).RowHeight = 30
$ws.Rows.Item(25).RowHeight = 30
$ws.Rows.Item(26).RowHeight = 30
